$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull data / push all data / mean calculation: update dSF (column F) values
$ws.Range("F2").Value = -4
$ws.Range("F4").Value = -6
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = 14
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 3
$ws.Range("F17").Value = 3
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 1
